$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the shared string "R40" (a label for rule R40).
# Update it to the text value "1" (kept as text, not a number).
$ws.Range("B11").Value = "1"
